$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "None"
$ws.Range("B3").Value  = "None"
$ws.Range("B4").Value  = "None"
$ws.Range("B7").Value  = "None"
$ws.Range("B8").Value  = "None"
$ws.Range("B9").Value  = "gNR-Bu"
$ws.Range("B10").Value = "None"
$ws.Range("B12").Value = "None"
$ws.Range("B13").Value = "None"
$ws.Range("B14").Value = "P-3O; P-4O; P-5O; P-6O"
$ws.Range("B15").Value = "None"
$ws.Range("B16").Value = "P3MEEET;"
$ws.Range("B18").Value = "None"
$ws.Range("B19").Value = "None"
$ws.Range("B20").Value = "p(g2T-T); Homo-gDPP"
$ws.Range("B21").Value = "None"
$ws.Range("B22").Value = "None"
$ws.Range("B23").Value = "PT-EG"
$ws.Range("B24").Value = "None"
$ws.Range("B25").Value = "None"
$ws.Range("B26").Value = "None"
$ws.Range("B27").Value = "None"
$ws.Range("B28").Value = "None"
$ws.Range("B29").Value = "oEG-substituted polythiophene"
$ws.Range("B30").Value = "PEDOT:PSS"
$ws.Range("B31").Value = "TIIP;"
$ws.Range("B32").Value = "None"
$ws.Range("B33").Value = "PT; NIDI"
$ws.Range("B34").Value = "None"
$ws.Range("B35").Value = "None"
$ws.Range("B36").Value = "None"
$ws.Range("B38").Value = "P3HT"
$ws.Range("B39").Value = "None"
$ws.Range("B40").Value = "None"
$ws.Range("B41").Value = "None"
$ws.Range("B43").Value = "None"
$ws.Range("B44").Value = "PEDOT:PSS"
$ws.Range("B46").Value = "PBBT-Me; BBL"
$ws.Range("B50").Value = "None"
$ws.Range("B55").Value = "None"
